$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row currently at position 30 (commessa 251651 / BIMEC 4) needs to move
# up to position 12 (its release date check should compare against
# "inizio lavorazione" rather than "fine lavorazione"), pushing rows 12-29
# down by one to 13-30.

# 1) Insert a blank row at 12 - this shifts the old rows 12..30 down to 13..31.
$ws.Rows("12:12").Insert()

# 2) The row that used to be row 30 is now row 31. Copy its values into the
#    newly-inserted row 12.
$v = $ws.Range("A31:S31").Value()
$ws.Range("A12:S12").Value = $v

# 3) Remove the now-duplicated row 31.
$ws.Rows("31:31").Delete()
